$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$range = $ws.Range("C2:C305")
$range.Value = 45178
